$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "67.607.84"
$ws.Range("E2").Value2 = "  +1.44%  "
$ws.Range("D3").Value2 = "3.964.61"
$ws.Range("E3").Value2 = "  +5.06%  "
$ws.Range("E4").Value2 = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "487.92"
$ws.Range("E5").Value2 = "  +10.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "149.48"
$ws.Range("E6").Value2 = "  +5.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.629"
$ws.Range("E7").Value2 = "  +1.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.998"
$ws.Range("E8").Value2 = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.735"
$ws.Range("E9").Value2 = "  +0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.168"
$ws.Range("E10").Value2 = "  +12.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.0000357"
$ws.Range("E11").Value2 = "  +16.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "43.46"
$ws.Range("E12").Value2 = "  +2.04%  "
$ws.Range("D13").Value2 = "4.580.44"
$ws.Range("E13").Value2 = "  +4.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "10.49"
$ws.Range("E14").Value2 = "  +1.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "15.10"
$ws.Range("E15").Value2 = "  +2.58%  "
$ws.Range("D16").Value2 = "3.974.32"
$ws.Range("E16").Value2 = "  +4.20%  "
$ws.Range("E17").Value2 = "  +0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "20.13"
$ws.Range("E18").Value2 = "  +1.91%  "
$ws.Range("E19").Value2 = "  +2.66%  "
$ws.Range("D20").Value2 = "67.711.47"
$ws.Range("E20").Value2 = "  +1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "435.50"
$ws.Range("E21").Value2 = "  +6.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "3.41"
$ws.Range("E22").Value2 = "  +5.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "14.64"
$ws.Range("E23").Value2 = "  +1.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "87.84"
$ws.Range("E24").Value2 = "  +3.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "3.73"
$ws.Range("E25").Value2 = "  +10.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "38.82"
$ws.Range("E26").Value2 = "  +6.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "10.17"
$ws.Range("E27").Value2 = "  +5.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "9.87"
$ws.Range("E28").Value2 = "  +1.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "717.99"
$ws.Range("E29").Value2 = "  -1.93%  "
$ws.Range("E30").Value2 = "  +0.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "13.44"
$ws.Range("E31").Value2 = "  -2.90%  "
$ws.Range("E32").Value2 = "  +3.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "42.46"
$ws.Range("E33").Value2 = "  -1.70%  "
$ws.Range("D34").Value2 = "0.0₃0851"
$ws.Range("E34").Value2 = "  +29.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "58.59"
$ws.Range("E35").Value2 = "  +4.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.151"
$ws.Range("E36").Value2 = "  -3.03%  "
$ws.Range("E37").Value2 = "  +0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "5.40"
$ws.Range("E38").Value2 = "  -2.15%  "
$ws.Range("E39").Value2 = "  +1.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "3.10"
$ws.Range("E40").Value2 = "  +7.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.145"
$ws.Range("E41").Value2 = "  +4.37%  "
$ws.Range("B42").Value2 = "ARBITRUM"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "2.25"
$ws.Range("E42").Value2 = "  +8.37%  "
$ws.Range("B43").Value2 = "TheGraph"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.338"
$ws.Range("E43").Value2 = "  +1.45%  "
$ws.Range("B44").Value2 = "WEMIXToken"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "2.84"
$ws.Range("E44").Value2 = "  +7.33%  "
$ws.Range("B45").Value2 = "FirstDigitalUSD"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.999"
$ws.Range("E45").Value2 = "  -0.27%  "
$ws.Range("E46").Value2 = "  +5.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "2.53"
$ws.Range("E47").Value2 = "  -4.69%  "
$ws.Range("B48").Value2 = "Monero"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "148.70"
$ws.Range("E48").Value2 = "  +4.24%  "
$ws.Range("B49").Value2 = "ApeXProtocol"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "3.22"
$ws.Range("E49").Value2 = "  -3.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "2.90"
$ws.Range("E50").Value2 = "  +3.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "25.71"
$ws.Range("E51").Value2 = "  +4.38%  "
